$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 100001530
$ws.Range("J17").Value = 100001530
$ws.Range("L17").Value = 300004590
$ws.Range("N17").Value = -300004926
$ws.Range("H96").Value = 232
$ws.Range("I96").Value = 238.28572
$ws.Range("K96").Value = 714.85716
$ws.Range("M96").Value = 658.14284
$ws.Range("H97").Value = 544.75
$ws.Range("J97").Value = 589.5
$ws.Range("L97").Value = 1768.5
$ws.Range("N97").Value = -2760.5
$ws.Range("H98").Value = 1588514.9
$ws.Range("I98").Value = 1852434
$ws.Range("K98").Value = 1852434
$ws.Range("M98").Value = -1850936
$ws.Range("H106").Value = 5467.25
$ws.Range("I106").Value = 5467.25
$ws.Range("K106").Value = 5467.25
$ws.Range("M106").Value = -4836.25
$ws.Range("H122").Value = 1588514.9
$ws.Range("I122").Value = 1852434
$ws.Range("K122").Value = 5557302
$ws.Range("M122").Value = -5554852
$ws.Range("H125").Value = 17963.666
$ws.Range("I125").Value = 1874
$ws.Range("J125").Value = 26008.5
$ws.Range("K125").Value = 16866
$ws.Range("L125").Value = 234076.5
$ws.Range("M125").Value = -14406
$ws.Range("N125").Value = -238996.5
$ws.Range("H137").Value = 1640.3334
$ws.Range("J137").Value = 2695
$ws.Range("L137").Value = 8085
$ws.Range("N137").Value = -13185
$ws.Range("H138").Value = 1978.5714
$ws.Range("I138").Value = 664.75
$ws.Range("K138").Value = 1994.25
$ws.Range("M138").Value = 3145.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 193.08333
$ws.Range("I4").Value = 226.55556
$ws.Range("J4").Value = 92.666664
$ws.Range("K4").Value = 226.55556
$ws.Range("L4").Value = 92.666664
$ws.Range("M4").Value = -110.55556
$ws.Range("N4").Value = -324.666664
$ws.Range("H32").Value = 9497.915000000001
$ws.Range("J32").Value = 18746.084
$ws.Range("L32").Value = 18746.084
$ws.Range("N32").Value = -19320.084
$ws.Range("H45").Value = 6574.222
$ws.Range("I45").Value = 4630.8
$ws.Range("J45").Value = 9003.5
$ws.Range("K45").Value = 4630.8
$ws.Range("L45").Value = 9003.5
$ws.Range("M45").Value = -4253.8
$ws.Range("N45").Value = -9757.5
$ws.Range("H74").Value = 2053.5151
$ws.Range("I74").Value = 1797
$ws.Range("K74").Value = 1797
$ws.Range("M74").Value = -923
$ws.Range("H77").Value = 2053.5151
$ws.Range("I77").Value = 1797
$ws.Range("K77").Value = 8985
$ws.Range("M77").Value = -4617
$ws.Range("H122").Value = 5770.852
$ws.Range("I122").Value = 4786.5454
$ws.Range("K122").Value = 14359.6362
$ws.Range("M122").Value = -11909.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 751.38464
$ws.Range("I107").Value = 661.8
$ws.Range("K107").Value = 661.8
$ws.Range("M107").Value = 1258.2
$ws.Range("H134").Value = 3819.7932
$ws.Range("I134").Value = 2298
$ws.Range("K134").Value = 6894
$ws.Range("M134").Value = -4359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1307.3636
$ws.Range("I16").Value = 375.6
$ws.Range("K16").Value = 375.6
$ws.Range("M16").Value = -88.60000000000002
$ws.Range("H31").Value = 4479.6665
$ws.Range("I31").Value = 2185.4443
$ws.Range("J31").Value = 5626.778
$ws.Range("K31").Value = 2185.4443
$ws.Range("L31").Value = 5626.778
$ws.Range("M31").Value = -1890.4443
$ws.Range("N31").Value = -6216.778
$ws.Range("H34").Value = 4479.6665
$ws.Range("I34").Value = 2185.4443
$ws.Range("J34").Value = 5626.778
$ws.Range("K34").Value = 2185.4443
$ws.Range("L34").Value = 5626.778
$ws.Range("M34").Value = -1983.4443
$ws.Range("N34").Value = -6030.778
$ws.Range("H113").Value = 1307.3636
$ws.Range("I113").Value = 375.6
$ws.Range("K113").Value = 375.6
$ws.Range("M113").Value = 1794.4
$ws.Range("H134").Value = 3515.3386
$ws.Range("I134").Value = 2506.5715
$ws.Range("K134").Value = 7519.7145
$ws.Range("M134").Value = -4984.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6729
$ws.Range("I113").Value = 10854.1
$ws.Range("J113").Value = 836
$ws.Range("K113").Value = 32562.3
$ws.Range("L113").Value = 2508
$ws.Range("M113").Value = -30392.3
$ws.Range("N113").Value = -6848
$ws.Range("H114").Value = 665.6111
$ws.Range("I114").Value = 291.7
$ws.Range("J114").Value = 1133
$ws.Range("K114").Value = 875.0999999999999
$ws.Range("L114").Value = 3399
$ws.Range("M114").Value = 2378.9
$ws.Range("N114").Value = -9907
$ws.Range("H117").Value = 389.41177
$ws.Range("J117").Value = 572.4
$ws.Range("L117").Value = 1717.2
$ws.Range("N117").Value = -8601.200000000001
$ws.Range("H129").Value = 2170.4119
$ws.Range("I129").Value = 1293.6
$ws.Range("K129").Value = 3880.8
$ws.Range("M129").Value = 1119.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3643.4583
$ws.Range("I132").Value = 3132.8667
$ws.Range("K132").Value = 9398.6001
$ws.Range("M132").Value = -6868.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4960.5386
$ws.Range("I7").Value = 2794.3333
$ws.Range("J7").Value = 5610.4
$ws.Range("K7").Value = 2794.3333
$ws.Range("L7").Value = 5610.4
$ws.Range("M7").Value = -2682.3333
$ws.Range("N7").Value = -5834.4
$ws.Range("H40").Value = 13762.143
$ws.Range("I40").Value = 35204
$ws.Range("K40").Value = 35204
$ws.Range("M40").Value = -35068
$ws.Range("H122").Value = 66671624
$ws.Range("I122").Value = 1000000000
$ws.Range("K122").Value = 3000000000
$ws.Range("M122").Value = -2999997550
$ws.Range("H126").Value = 4960.5386
$ws.Range("I126").Value = 2794.3333
$ws.Range("J126").Value = 5610.4
$ws.Range("K126").Value = 8382.999899999999
$ws.Range("L126").Value = 16831.2
$ws.Range("M126").Value = -5912.999899999999
$ws.Range("N126").Value = -21771.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 19986.5
$ws.Range("J76").Value = 19986.5
$ws.Range("L76").Value = 19986.5
$ws.Range("N76").Value = -20616.5
$ws.Range("H79").Value = 19986.5
$ws.Range("J79").Value = 19986.5
$ws.Range("L79").Value = 19986.5
$ws.Range("N79").Value = -22170.5
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 8945.071
$ws.Range("I81").Value = 11182.7
$ws.Range("K81").Value = 22365.4
$ws.Range("M81").Value = -21304.4
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 8945.071
$ws.Range("I84").Value = 11182.7
$ws.Range("K84").Value = 111827
$ws.Range("M84").Value = -106523
$ws.Range("H126").Value = 2695.0908
$ws.Range("I126").Value = 2825
$ws.Range("J126").Value = 2110.5
$ws.Range("K126").Value = 8475
$ws.Range("L126").Value = 6331.5
$ws.Range("M126").Value = -6005
$ws.Range("N126").Value = -11271.5
